$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 107799.9
$ws.Range("I62").Value = 206059.8
$ws.Range("J62").Value = 9540
$ws.Range("K62").Value = 206059.8
$ws.Range("L62").Value = 9540
$ws.Range("M62").Value = -205435.8
$ws.Range("N62").Value = -10788
$ws.Range("H65").Value = 107799.9
$ws.Range("I65").Value = 206059.8
$ws.Range("J65").Value = 9540
$ws.Range("K65").Value = 1030299
$ws.Range("L65").Value = 47700
$ws.Range("M65").Value = -1027179
$ws.Range("N65").Value = -53940
$ws.Range("H121").Value = 800
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H138").Value = 4090.15
$ws.Range("I138").Value = 902.64703
$ws.Range("J138").Value = 4743.012
$ws.Range("K138").Value = 2707.94109
$ws.Range("L138").Value = 14229.036
$ws.Range("M138").Value = 2432.05891
$ws.Range("N138").Value = -24509.036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2795
$ws.Range("I2").Value = 3827.75
$ws.Range("J2").Value = 2106.5
$ws.Range("K2").Value = 3827.75
$ws.Range("L2").Value = 2106.5
$ws.Range("M2").Value = -3714.75
$ws.Range("N2").Value = -2332.5
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30748
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164
$ws.Range("H102").Value = 3300
$ws.Range("I102").Value = 3750
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 3750
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -2128
$ws.Range("N102").Value = -4744
$ws.Range("H116").Value = 2795
$ws.Range("I116").Value = 3827.75
$ws.Range("J116").Value = 2106.5
$ws.Range("K116").Value = 3827.75
$ws.Range("L116").Value = 2106.5
$ws.Range("M116").Value = -1533.75
$ws.Range("N116").Value = -6694.5
$ws.Range("H132").Value = 3332.8333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3332.8333
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 9998.499899999999
$ws.Range("N132").Value = -15058.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2795
$ws.Range("I3").Value = 3827.75
$ws.Range("J3").Value = 2106.5
$ws.Range("K3").Value = 3827.75
$ws.Range("L3").Value = 2106.5
$ws.Range("M3").Value = -3713.75
$ws.Range("N3").Value = -2334.5
$ws.Range("H132").Value = 161444.44
$ws.Range("J132").Value = 161444.44
$ws.Range("L132").Value = 161444.44
$ws.Range("N132").Value = -171564.44
$ws.Range("H134").Value = 121308.12
$ws.Range("I134").Value = 4762.4
$ws.Range("J134").Value = 287802
$ws.Range("K134").Value = 14287.2
$ws.Range("L134").Value = 863406
$ws.Range("M134").Value = -11752.2
$ws.Range("N134").Value = -868476
$ws.Range("H135").Value = 70779.5
$ws.Range("J135").Value = 70779.5
$ws.Range("L135").Value = 70779.5
$ws.Range("N135").Value = -80919.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2263.75
$ws.Range("I99").Value = 1754.3478
$ws.Range("J99").Value = 2952.9412
$ws.Range("K99").Value = 1754.3478
$ws.Range("L99").Value = 2952.9412
$ws.Range("M99").Value = -256.3478
$ws.Range("N99").Value = -5948.9412
$ws.Range("H100").Value = 35057.5
$ws.Range("J100").Value = 35057.5
$ws.Range("L100").Value = 35057.5
$ws.Range("N100").Value = -37221.5
$ws.Range("H126").Value = 2263.75
$ws.Range("I126").Value = 1754.3478
$ws.Range("J126").Value = 2952.9412
$ws.Range("K126").Value = 5263.0434
$ws.Range("L126").Value = 8858.8236
$ws.Range("M126").Value = -2793.0434
$ws.Range("N126").Value = -13798.8236
$ws.Range("H132").Value = 2619.7273
$ws.Range("I132").Value = 2029.7333
$ws.Range("J132").Value = 3884
$ws.Range("K132").Value = 6089.199900000001
$ws.Range("L132").Value = 11652
$ws.Range("M132").Value = -3559.199900000001
$ws.Range("N132").Value = -16712
$ws.Range("H134").Value = 2537.75
$ws.Range("I134").Value = 2052.6785
$ws.Range("K134").Value = 6158.0355
$ws.Range("M134").Value = -3623.0355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 726.5
$ws.Range("I76").Value = 726.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2179.5
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -1796.5
$ws.Range("H79").Value = 726.5
$ws.Range("I79").Value = 726.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2179.5
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -853.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 405
$ws.Range("I31").Value = 405
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 405
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -113
$ws.Range("H37").Value = 405
$ws.Range("I37").Value = 405
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 405
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -128
$ws.Range("H70").Value = 4675.5835
$ws.Range("I70").Value = 4014
$ws.Range("J70").Value = 5601.8
$ws.Range("K70").Value = 4014
$ws.Range("L70").Value = 5601.8
$ws.Range("M70").Value = -3744
$ws.Range("N70").Value = -6141.8
$ws.Range("H73").Value = 4675.5835
$ws.Range("I73").Value = 4014
$ws.Range("J73").Value = 5601.8
$ws.Range("K73").Value = 4014
$ws.Range("L73").Value = 5601.8
$ws.Range("M73").Value = -3078
$ws.Range("N73").Value = -7473.8
$ws.Range("H132").Value = 3663.077
$ws.Range("I132").Value = 3337.3333
$ws.Range("J132").Value = 3942.2856
$ws.Range("K132").Value = 10011.9999
$ws.Range("L132").Value = 11826.8568
$ws.Range("M132").Value = -7481.999899999999
$ws.Range("N132").Value = -16886.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 18520786
$ws.Range("I100").Value = 22224502
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 22224502
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -22223961
$ws.Range("N100").Value = -3282
$ws.Range("H132").Value = 1923.28
$ws.Range("J132").Value = 3200.375
$ws.Range("L132").Value = 9601.125
$ws.Range("M132").Value = -1436.8823
$ws.Range("N132").Value = -14661.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3450
$ws.Range("I96").Value = 3450
$ws.Range("K96").Value = 3450
$ws.Range("M96").Value = -2077
$ws.Range("H135").Value = 78373.5
$ws.Range("J135").Value = 78373.5
$ws.Range("L135").Value = 78373.5
$ws.Range("N135").Value = -88513.5
